$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new columns C and F by shifting header/data ---
# First, copy style from A1 (bold header style) into the new header cells C1, F1 and G1
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("A1").Copy($ws.Range("G1"))

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "loss"
$ws.Range("B1").Value = "binary_io_u_3"
$ws.Range("C1").Value = "dice_metric"
$ws.Range("D1").Value = "val_loss"
$ws.Range("E1").Value = "val_binary_io_u_3"
$ws.Range("F1").Value = "val_dice_metric"
$ws.Range("G1").Value = "epoch"

# --- Update data rows (rows 2-20) ---
$ws.Range("A2").Value = 1.065317153930664
$ws.Range("B2").Value = 0.4640493988990784
$ws.Range("C2").Value = 0.3548171520233154
$ws.Range("D2").Value = 1.271228075027466
$ws.Range("E2").Value = 0.4574995040893555
$ws.Range("F2").Value = 0.2631695866584778
$ws.Range("G2").Value = 1

$ws.Range("A3").Value = 0.8141606450080872
$ws.Range("B3").Value = 0.5656570196151733
$ws.Range("C3").Value = 0.4668542444705963
$ws.Range("D3").Value = 1.176835894584656
$ws.Range("E3").Value = 0.4574995040893555
$ws.Range("F3").Value = 0.2735696732997894
$ws.Range("G3").Value = 2

$ws.Range("A4").Value = 0.6438578963279724
$ws.Range("B4").Value = 0.6217673420906067
$ws.Range("C4").Value = 0.5723934769630432
$ws.Range("D4").Value = 1.089979887008667
$ws.Range("E4").Value = 0.4574995040893555
$ws.Range("F4").Value = 0.2712994515895844
$ws.Range("G4").Value = 3

$ws.Range("A5").Value = 0.4910151362419128
$ws.Range("B5").Value = 0.6713499426841736
$ws.Range("C5").Value = 0.6832207441329956
$ws.Range("D5").Value = 1.056345224380493
$ws.Range("E5").Value = 0.4574995040893555
$ws.Range("F5").Value = 0.24927918612957
$ws.Range("G5").Value = 4

$ws.Range("A6").Value = 0.4103506505489349
$ws.Range("B6").Value = 0.6957248449325562
$ws.Range("C6").Value = 0.7400398254394531
$ws.Range("D6").Value = 1.055898308753967
$ws.Range("E6").Value = 0.4574995040893555
$ws.Range("F6").Value = 0.2041231691837311
$ws.Range("G6").Value = 5

$ws.Range("A7").Value = 0.3730241358280182
$ws.Range("B7").Value = 0.7027483582496643
$ws.Range("C7").Value = 0.7659919857978821
$ws.Range("D7").Value = 1.081430196762085
$ws.Range("E7").Value = 0.4574995040893555
$ws.Range("F7").Value = 0.1559834629297256
$ws.Range("G7").Value = 6

$ws.Range("A8").Value = 0.3533163666725159
$ws.Range("B8").Value = 0.7108031511306763
$ws.Range("C8").Value = 0.7779204845428467
$ws.Range("D8").Value = 1.115687012672424
$ws.Range("E8").Value = 0.4574995040893555
$ws.Range("F8").Value = 0.1152432933449745
$ws.Range("G8").Value = 7

$ws.Range("A9").Value = 0.3328442871570587
$ws.Range("B9").Value = 0.7127156853675842
$ws.Range("C9").Value = 0.7922288179397583
$ws.Range("D9").Value = 1.163525700569153
$ws.Range("E9").Value = 0.4574995040893555
$ws.Range("F9").Value = 0.07983924448490143
$ws.Range("G9").Value = 8

$ws.Range("A10").Value = 0.3264735341072083
$ws.Range("B10").Value = 0.7162885665893555
$ws.Range("C10").Value = 0.7952737808227539
$ws.Range("D10").Value = 1.191913723945618
$ws.Range("E10").Value = 0.4574995040893555
$ws.Range("F10").Value = 0.06144086644053459
$ws.Range("G10").Value = 9

$ws.Range("A11").Value = 0.3187452554702759
$ws.Range("B11").Value = 0.7206717729568481
$ws.Range("C11").Value = 0.7989698648452759
$ws.Range("D11").Value = 1.229954123497009
$ws.Range("E11").Value = 0.4574995040893555
$ws.Range("F11").Value = 0.04168339446187019
$ws.Range("G11").Value = 10

$ws.Range("A12").Value = 0.3118893504142761
$ws.Range("B12").Value = 0.7166717052459717
$ws.Range("C12").Value = 0.8046832084655762
$ws.Range("D12").Value = 1.257521867752075
$ws.Range("E12").Value = 0.4574995040893555
$ws.Range("F12").Value = 0.03299479931592941
$ws.Range("G12").Value = 11

$ws.Range("A13").Value = 0.3074690699577332
$ws.Range("B13").Value = 0.7213510274887085
$ws.Range("C13").Value = 0.8062549829483032
$ws.Range("D13").Value = 1.289356708526611
$ws.Range("E13").Value = 0.4574995040893555
$ws.Range("F13").Value = 0.02305548079311848
$ws.Range("G13").Value = 12

$ws.Range("A14").Value = 0.3030803203582764
$ws.Range("B14").Value = 0.7233078479766846
$ws.Range("C14").Value = 0.8082906603813171
$ws.Range("D14").Value = 1.319879174232483
$ws.Range("E14").Value = 0.4574995040893555
$ws.Range("F14").Value = 0.01624204032123089
$ws.Range("G14").Value = 13

$ws.Range("A15").Value = 0.3050404191017151
$ws.Range("B15").Value = 0.7186112403869629
$ws.Range("C15").Value = 0.8076108694076538
$ws.Range("D15").Value = 1.3380286693573
$ws.Range("E15").Value = 0.4574995040893555
$ws.Range("F15").Value = 0.01414340455085039
$ws.Range("G15").Value = 14

$ws.Range("A16").Value = 0.2937130630016327
$ws.Range("B16").Value = 0.7252293825149536
$ws.Range("C16").Value = 0.8148374557495117
$ws.Range("D16").Value = 1.343879342079163
$ws.Range("E16").Value = 0.4574995040893555
$ws.Range("F16").Value = 0.01307932287454605
$ws.Range("G16").Value = 15

$ws.Range("A17").Value = 0.2910315692424774
$ws.Range("B17").Value = 0.7261908054351807
$ws.Range("C17").Value = 0.8162949681282043
$ws.Range("D17").Value = 1.366653919219971
$ws.Range("E17").Value = 0.4574995040893555
$ws.Range("F17").Value = 0.01036407984793186
$ws.Range("G17").Value = 16

$ws.Range("A18").Value = 0.29124516248703
$ws.Range("B18").Value = 0.7235106825828552
$ws.Range("C18").Value = 0.8162726163864136
$ws.Range("D18").Value = 1.386469602584839
$ws.Range("E18").Value = 0.4574995040893555
$ws.Range("F18").Value = 0.008140009827911854
$ws.Range("G18").Value = 17

$ws.Range("A19").Value = 0.2891172766685486
$ws.Range("B19").Value = 0.7251983880996704
$ws.Range("C19").Value = 0.8174265027046204
$ws.Range("D19").Value = 1.392734289169312
$ws.Range("E19").Value = 0.4574995040893555
$ws.Range("F19").Value = 0.007918978109955788
$ws.Range("G19").Value = 18

$ws.Range("A20").Value = 0.2851462960243225
$ws.Range("B20").Value = 0.727731466293335
$ws.Range("C20").Value = 0.8196375966072083
$ws.Range("D20").Value = 1.386369943618774
$ws.Range("E20").Value = 0.457520067691803
$ws.Range("F20").Value = 0.008769180625677109
$ws.Range("G20").Value = 19

